$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Have" (C column) quantities for rows 2-13 (skip row 6)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1

# Update E column formulas to clamp at 0 using MAX()
$ws.Range("E2").Formula = "=MAX(B2-(C2+D2),0)"
$ws.Range("E3:E32").Formula = "=MAX(B3-(C3+D3),0)"

# Update selected cell
$ws.Range("C11").Select()
